$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16: change the plenary table's table style.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CE7E20AD-3B8D-4FDA-B002-E014FF2973ED}")
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the deck's main theme palette ("Integral") for the notes-master's
#    "Office Theme" palette (the two theme parts share an identical font /
#    format scheme, so only the 12 colour-scheme slots actually differ).
# ---------------------------------------------------------------------------
function Set-ThemeColor($scheme, $idx, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Colors($idx).RGB = $r + ($g * 256) + ($b * 65536)
}

$tcs = $slide.ThemeColorScheme

Set-ThemeColor $tcs 1  "000000"   # dk1
Set-ThemeColor $tcs 2  "FFFFFF"   # lt1
Set-ThemeColor $tcs 3  "44546A"   # dk2
Set-ThemeColor $tcs 4  "E7E6E6"   # lt2
Set-ThemeColor $tcs 5  "5B9BD5"   # accent1
Set-ThemeColor $tcs 6  "ED7D31"   # accent2
Set-ThemeColor $tcs 7  "A5A5A5"   # accent3
Set-ThemeColor $tcs 8  "FFC000"   # accent4
Set-ThemeColor $tcs 9  "4472C4"   # accent5
Set-ThemeColor $tcs 10 "70AD47"   # accent6
Set-ThemeColor $tcs 11 "0563C1"   # hlink
Set-ThemeColor $tcs 12 "954F72"   # folHlink
